# Trade #68 closed at 2026-02-17 08:52:13 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the newly closed trade #68.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value2 = 1200.26   # Current Capital
$summary.Range("B6").Value2 = 68        # Total Trades
$summary.Range("B9").Value2 = 41.18     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value2 = 100.26     # Capital
$status.Range("D4").Value2 = 68         # Trades
$status.Range("F4").Value2 = 0.26       # P&L %
$status.Range("G4").Value2 = 41.18      # Win Rate %

# ---------------------------------------------------------------------
# Helper to append the new trade row (#68) to a trades sheet.
# Date/time columns must stay plain text (not Excel date serials), so the
# number format is forced to Text before the values are written.
# ---------------------------------------------------------------------
function Add-Trade68($ws) {
    $row = 69

    $ws.Cells.Item($row, 1).Value2 = 68                # Trade #

    # Force text format on the date column only - Excel would otherwise
    # coerce a "yyyy-mm-dd" literal into a date serial number. Plain
    # "hh:mm:ss" strings are not re-interpreted, so column C needs no
    # special handling.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value2 = "2026-02-17"       # Date

    $ws.Cells.Item($row, 3).Value2 = "08:52:07"         # Time

    $ws.Cells.Item($row, 4).Value2 = "MarketMaking"     # Strategy
    $ws.Cells.Item($row, 5).Value2 = "UP"               # Side
    $ws.Cells.Item($row, 6).Value2 = 0.41               # Entry Price
    $ws.Cells.Item($row, 7).Value2 = 0.405941           # Exit Price
    $ws.Cells.Item($row, 8).Value2 = "CLOSED"           # Status
    $ws.Cells.Item($row, 9).Value2 = -0.9901            # P&L %
    $ws.Cells.Item($row, 10).Value2 = -0                # P&L $
    $ws.Cells.Item($row, 11).Value2 = 100.26            # Capital After
    $ws.Cells.Item($row, 12).Value2 = 0                 # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value2 = 0                 # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value2 = 0.6               # Confidence
    $ws.Cells.Item($row, 15).Value2 = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value2 = "early_exit"      # Exit Reason
    $ws.Cells.Item($row, 17).Value2 = 0.13              # Duration (min)
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade68 $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade68 $marketMaking
